# daily auto push: 2025-10-10 18:38 UTC
# Append the next day's row (row 92) to the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 92

# A92 holds a date-like string ("2025/10/11"); format the cell as Text
# first so Excel stores it as a literal string instead of auto-converting
# it to a date serial number, then drop the formatting again so the cell
# ends up with the same (default) style as its neighbours.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2025/10/11"
$ws.Range("A$row").ClearFormats()

$ws.Range("B$row").Value = "土"
$ws.Range("C$row").Value = 0
$ws.Range("D$row").Value = 201
